$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 0. Snapshot every existing comment (row, column, text) before touching
#    anything. This engine's Rows.Delete() does not re-anchor comments the
#    way Excel does, so the move below is done by hand: capture -> delete
#    row -> remove all old comments -> re-add at the shifted row.
# ---------------------------------------------------------------------------
$oldComments = $ws.Comments
$commentCount = $oldComments.Count
$commentRows = @()
$commentCols = @()
$commentTexts = @()
for ($i = 1; $i -le $commentCount; $i++) {
    $cm = $oldComments.Item($i)
    $p = $cm.Parent
    $commentRows += $p.Row
    $commentCols += $p.Column
    $commentTexts += $cm.Text()
}

# ---------------------------------------------------------------------------
# 1. M07 list: drop the "Produceren van een SBoM" bullet and renumber the
#    final bullet from "10." to "9.". At this point (before the row delete
#    below) that final bullet still lives in B52.
# ---------------------------------------------------------------------------
$ws.Range("B52").Value2 = "9. Oplevering van het totale product, dus inclusief alle deliverables, in de vorm zoals bruikbaar voor en afgesproken met de opdrachtgevende organisatie."

# Remove the now-obsolete "9. Produceren van een ""software bill of materials"" (SBoM)," row (B51).
# Deleting the whole row shifts every row below it up by one (matching the rest of the diff).
$ws.Rows("51").Delete()

# ---------------------------------------------------------------------------
# 2. M16 list items (now B53:B69 after the row delete above): strip the
#    trailing ": <tool>" part and capitalise the first letter of each task,
#    since the tool names move into the comment's new table instead.
# ---------------------------------------------------------------------------
$ws.Range("B53").Value2 = "1. Product en sprint backlog management en agile werken"
$ws.Range("B54").Value2 = "2. Inrichten en uitvoeren van een continuous delivery pipeline"
$ws.Range("B55").Value2 = "3. Monitoren van de kwaliteit van broncode"
$ws.Range("B56").Value2 = "4. Versiebeheer van op te leveren producten"
$ws.Range("B57").Value2 = "5. Release van software"
$ws.Range("B58").Value2 = "6. Maken van testrapportages"
$ws.Range("B59").Value2 = "7. Maken van kwaliteitsrapportages"
$ws.Range("B60").Value2 = "8. Controleren op aanwezigheid van bekende kwetsbaarheden in externe software"
$ws.Range("B61").Value2 = "9. Statische controle van de software op aanwezigheid van kwetsbare constructies"
$ws.Range("B62").Value2 = "10. Dynamische controle van de software op aanwezigheid van kwetsbare constructies"
$ws.Range("B63").Value2 = "11. Controleren van container images op aanwezigheid van bekende kwetsbaarheden"
$ws.Range("B64").Value2 = "12. Testen van performance en schaalbaarheid"
$ws.Range("B65").Value2 = "13. Testen op toegankelijkheid van de applicatie"
$ws.Range("B66").Value2 = "14. Produceren van een ""software bill of materials"" (SBoM)"
$ws.Range("B67").Value2 = "15. Opslaan van artifacten"
$ws.Range("B68").Value2 = "16. Registratie van incidenten bij gebruik en beheer"
$ws.Range("B69").Value2 = "17. Bij het uitvoeren van operationeel beheer; uitrollen van de software in de productieomgeving"

# ---------------------------------------------------------------------------
# 3. New text for the three comments whose content changes. Everything else
#    keeps its captured text verbatim.
# ---------------------------------------------------------------------------
$m07Comment = "M07: Het project gebruikt een continuous delivery pipeline om het product te bouwen, testen en op te leveren`n`nEr is een geautomatiseerde continuous delivery pipeline die aantoonbaar correct werkt en de software bouwt, installeert in de testomgevingen, test op functionele en niet-functionele eigenschappen en oplevert, al dan niet inclusief installatie in de productieomgeving.`n`nDe geautomatiseerde continuous delivery pipeline voert ten minste de volgende activiteiten uit:`n`n1. Bouw van de software,`n2. Unit tests,`n3. Regressietests,`n4. Beveiligingstests,`n5. Performancetests,`n6. Toegankelijkheidstests,`n7. Broncodekwaliteitscontroles,`n8. Installatie van de software in test, acceptatie en/of productieomgevingen,`n9. Oplevering van het totale product, dus inclusief alle deliverables, in de vorm zoals bruikbaar voor en afgesproken met de opdrachtgevende organisatie.`n`nPerformance- en beveiligingstests op de software zijn ook onderdeel van de continuous delivery pipeline, maar vanwege doorlooptijden en licenties is dat niet altijd haalbaar; in dat geval vinden de performance- en beveiligingstests zo veel mogelijk, en bij voorkeur dagelijks, plaats. Performance- en beveiligingstests op de software vinden plaats in de testomgeving van het project. Als ICTU verantwoordelijk is voor het operationeel beheer laat ICTU de performance- en beveiligingstesten op de software (ook) uitvoeren in een productie-like omgeving.`n`nNiet alle testen en controles kunnen altijd geautomatiseerd worden uitgevoerd. Denk aan kwaliteitscontroles op architectuurbeslissingen of het testen van toegankelijkheidseisen. Waar mogelijk wordt wel een zo groot mogelijk deel van de testen en controles geautomatiseerd en als onderdeel van de pipeline uitgevoerd.`n`nRationale`n`nSoftware incrementeel opleveren vereist dat de software frequent gebouwd, getest en opgeleverd kan worden. Om dit efficiënt en foutvrij te doen, dient het proces van bouwen, testen en opleveren geautomatiseerd te zijn; een continuous delivery pipeline faciliteert dit.`n"

$m16Comment = "M16: Het project gebruikt tools voor vastgestelde taken`n`nVoor vastgestelde taken bij het ontwikkelen, onderhouden en operationeel beheren van software, stelt ICTU het gebruik van tools verplicht. ICTU adviseert per taak specifieke tools en ondersteunt projecten bij het gebruik daarvan.`n`nICTU adviseert en ondersteunt voor de hieronder genoemde taken specifieke tools. Projecten gebruiken deze tools, of gelijkwaardige alternatieven.`n`nActiviteit                                                                                   Tools                                                                                    `nProduct en sprint backlog management en agile werken                                         Azure DevOps of Jira                                                                     `nInrichten en uitvoeren van een continuous delivery pipeline                                  Jenkins, GitLab CI/CD (Continuous Integration, Delivery, and Deployment) of Azure DevOps `nMonitoren van de kwaliteit van broncode                                                      SonarQube                                                                                `nVersiebeheer van op te leveren producten                                                     GitLab of Azure DevOps                                                                   `nRelease van software                                                                         Releaseserver in het ontwikkelplatform                                                   `nMaken van testrapportages                                                                    JUnit, Robot Framework, TestNG, of hiermee compatible tools                              `nMaken van kwaliteitsrapportages                                                              Quality-time                                                                             `nControleren op aanwezigheid van bekende kwetsbaarheden in externe software                   OWASP (Open Web Application Security Project) Dependency-Check en/of Dependency-Track    `nStatische controle van de software op aanwezigheid van kwetsbare constructies                SonarQube                                                                                `nDynamische controle van de software op aanwezigheid van kwetsbare constructies               ZAP (Zed Attack Proxy) by Checkmarx                                                      `nControleren van container images op aanwezigheid van bekende kwetsbaarheden                  Trivy                                                                                    `nTesten van performance en schaalbaarheid                                                     JMeter en Performancetestrunner                                                          `nTesten op toegankelijkheid van de applicatie                                                 Axe                                                                                      `nProduceren van een ""software bill of materials"" (SBoM)                                       Tools die een SBoM in CycloneDX-formaat (zie https://cyclonedx.org) genereren            `nOpslaan van artifacten                                                                       Nexus of Harbor                                                                          `nRegistratie van incidenten bij gebruik en beheer                                             Jira                                                                                     `nBij het uitvoeren van operationeel beheer; uitrollen van de software in de productieomgeving Ansible                                                                                  `nN.B. Onder het ondersteunen van ""agile werken"" vallen het opvoeren van eisen, het opvoeren van logische testgevallen, het koppelen van logische testgevallen aan eisen, het bijhouden van een werkvoorraad, het plannen van iteraties en het toewijzen van eisen aan iteraties. De 'eisen' worden, conform Scrumterminologie, geregistreerd als epics en/of user stories, de werkvoorraad als product backlog en de iteraties als sprints. Het toewijzen van eisen aan iteraties gebeurt via de sprint backlog.`n`nRationale`n`nProjecten hebben een redelijke vrijheid bij het kiezen en gebruiken van tools, maar voor een aantal taken is het gebruik verplicht gesteld. Deze tools zijn nodig voor een efficiënte uitvoering van de Kwaliteitsaanpak. Uniform gebruik van deze tools maakt het mogelijk koppeling tussen die tools voor alle projecten te standaardiseren; daarnaast bevordert het de uitwisselbaarheid van medewerkers en neemt het risico op het gebruik van onvolwassen tools af. Tot slot is het gebruik in een aantal gevallen, ten behoeve van informatiebeveiliging bij de overheid, verplicht.`n"

$m18Comment = "M18: ICTU biedt ondersteuning voor verplicht gestelde tools`n`nICTU zorgt voor technische en functionele ondersteuning aan projecten bij het gebruik van alle verplichte tools.`n`nICTU zorgt voor ondersteuning van de in M16: Het project gebruikt tools voor vastgestelde taken verplicht gestelde tools. Een team van specialisten met kennis, ervaring en capaciteit is beschikbaar voor ondersteuning aan projecten. Projecten zijn verantwoordelijk voor de correcte werking van de pipeline.`n`nBij de selectie van tools ter ondersteuning van de projectuitvoering geeft ICTU de voorkeur aan open source tools. Ook tools die ICTU zelf ontwikkelt ter ondersteuning van softwareontwikkelprojecten worden bij voorkeur open source beschikbaar gesteld.`n`nRationale`n`nDe keuze om het gebruik van een aantal tools verplicht te stellen (M16: Het project gebruikt tools voor vastgestelde taken) volgt uit de belangrijke rol die die tools spelen in de ontwikkelstraat en in Quality-time, het kwaliteitssysteem van ICTU. Met de verplichting komt ook een verantwoordelijkheid: om projecten in staat te stellen snel en effectief met deze tools te werken, moeten die projecten ondersteund worden.`n`nDe verplicht gestelde tools zijn beperkt in aantal, bewezen en gangbaar; veel medewerkers zullen deze tools al kennen.`n`nDe voorkeur voor open source tools is conform de rationale uit NORA (Nederlandse Overheid Referentiearchitectuur) voor het gebruik van open source tools, zoals beschreven in NORA v3.0 drijfveer ""Beleid open standaarden"". De voorkeur voor het open source beschikbaar stellen van eigen ontwikkelde tools is conform de ""Beleidsbrief vrijgeven van de broncode van overheidssoftware"" van de staatssecretaris van Binnenlandse Zaken en Koninkrijksrelaties, 17 april 2020.`n"

# ---------------------------------------------------------------------------
# 4. Remove every old comment shape (their anchors are now stale for
#    anything that was on row 52 or below), then re-add each one on its
#    (possibly shifted) row, substituting the new text where it changed.
# ---------------------------------------------------------------------------
$liveComments = $ws.Comments
for ($i = $liveComments.Count; $i -ge 1; $i--) {
    $liveComments.Item($i).Delete()
}

for ($i = 0; $i -lt $commentCount; $i++) {
    $r = $commentRows[$i]
    $c = $commentCols[$i]
    if ($r -ge 52) {
        $r = $r - 1
    }
    $text = $commentTexts[$i]
    if ($r -eq 42 -and $c -eq 2) {
        $text = $m07Comment
    } elseif ($r -eq 52 -and $c -eq 2) {
        $text = $m16Comment
    } elseif ($r -eq 104 -and $c -eq 2) {
        $text = $m18Comment
    }
    $cell = $ws.Cells.Item($r, $c)
    $cell.AddComment($text)
}
